# level_main_07-05_end: normalize curly quotes / NBSPs to straight quotes / spaces
# in the en_US (col C) and ko_KR (col D) translation columns.
# Commit: "update on 20210731 画中人"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- en_US column (C) ---

# "Bloodied. Torn. Pierced.<NBSP>" -> trailing NBSP becomes a regular space
$ws.Range("C2").Value = "Bloodied. Torn. Pierced. `n"

# "[Decision(options=""How old is she?<NBSP>"", ...)]" -> NBSP becomes a regular space
$ws.Range("C14").Value = "[Decision(options=""How old is she? "", values=""1"")]`n"

# "[name=""Kal'tsit""]  Fourteen.<NBSP>" -> NBSP becomes a regular space
$ws.Range("C15").Value = "[name=""Kal'tsit""]  Fourteen. `n"

# curly double quotes “...” -> straight single quotes '...'
$ws.Range("C40").Value = "[name=""Rosmontis""]  When Amiya or I stand on the battlefield... who looks at us and thinks 'children?'`n"
$ws.Range("C41").Value = "[name=""Rosmontis""]  We’re 'monsters,' aren’t we, Doctor?`n"

$ws.Range("C118").Value = "[name=""Kal'tsit""]  'When you deprive a human being of her sentience, what is left? What is created?'`n"
$ws.Range("C119").Value = "[name=""Kal'tsit""]  'When that creation deprives us of our lives, who is at fault?'`n"

# "...more.<NBSP>" -> NBSP becomes a regular space
$ws.Range("C122").Value = "[name=""Kal'tsit""]  As we move on to the next battle, I may be able to teach you more. `n"

$ws.Range("C148").Value = "[name=""Kal'tsit""]  ...'The last Wendigo will die at the hands of the Lord of Fiends?'`n"
$ws.Range("C150").Value = "[name=""Kal'tsit""]  The original phrasing of the prophecy was 'The son of Hor-Tekrz, traitor to Sarkaz and disreputable end of the bloodline, will be executed by the Lord of Sarkaz.'`n"
$ws.Range("C151").Value = "[name=""Kal'tsit""]  Some prophecy. Whether or not there’s any extension of their 'bloodline,' there are still many Wendigos living their lives in Columbia and Ursus.`n"
$ws.Range("C152").Value = "[name=""Kal'tsit""]  ...That is, if a few dozen can be called 'many.'`n"
$ws.Range("C155").Value = "[name=""Kal'tsit""]  Are you getting old? Since when do you believe in 'prophecies' and Sarkaz witchcraft?`n"
$ws.Range("C176").Value = "[name=""Kal'tsit""]  Some time ago, an Infected clinic in Chernobog called 'Azazel' traded information with us. I learned some more information about the Wendigo in Reunion's service.`n"

# outer curly double quotes -> straight single quotes; inner curly single quotes kept as-is
$ws.Range("C184").Value = "[name=""Kal'tsit""]  'Do not, under any circumstances, even at Rhodes Island, speak of ‘the Lord of Fiends.’'`n"

# --- ko_KR column (D) ---
# 흡혈귀 ("vampire", Sino-Korean) -> 뱀파이어 ("vampire", loanword)
$ws.Range("D145").Value = "[name=""켈시""]  카우투스의 생리적 변화는 뱀파이어와는 차원이 달라. 우리가 카즈델을 떠난 지 벌써 삼 년이라고.`n"
$ws.Range("D164").Value = "[name=""와파린""]  나는 뱀파이어고, 그 늙은이는 웬디고라고 해도 말이야.`n"
